# Review_406.docx edit: swap in the "Addition Is All You Need" (L-Mul) review
# content in place of the "Understanding Visual Feature Reliance" review,
# trimming the document down to the new, shorter set of paragraphs.

$d = $word.ActiveDocument
$br = [char]11   # manual line break (w:br) marker used by Word's text model

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find.Execute failed to locate: $old"
    }
}

# --- Paragraph 1: title / date line ---------------------------------------
Replace-Text "המאמר היומי של מייק - 25.02.25" "המאמר היומי של מייק - 23.02.25"
Replace-Text "Understanding Visual Feature Reliance through the Lens of Complexity" `
             "Addition Is All You Need: For Energy-Efficient Language Models"

# --- Paragraph 2: intro ------------------------------------------------------
Replace-Text "המאמר שאני סוקר היום מציג מחקר יוצא דופן, נדיר ומעניין על מורכבות פיצ'רים המופקים על ידי מודלים דיפ (אין RAG, סוכנים ו-LLMs שם :). מאמר זה קשור הדוקות לרעיון של צוואר הבקבוק של המידע ברשתות עצביות עמוקות, שטבע נפתלי תשבי." `
             ("מבוא:" + $br + "המאמר מציג גישה אלגנטית אך רדיקלית לשיפור היעילות של רשתות נוירונים, רלוונטית במיוחד לשיפור ביצועים של LLMs. המחברים מציעים חלופה למכפלות נקודה צפה(floating point) מסורתיות (Linear-Complexity Multiplication(L-Mul, אשר מקרב פעולות עם נקודה צפה על ידי חיבורי מספרים שלמים. הטענה המרכזית היא ש-L-Mul מפחית משמעותית את המורכבות החישובית ואת צריכת האנרגיה, תוך שמירה על ביצועי מודל כמעט זהים.")

# --- Paragraph 3: motivation -------------------------------------------------
Replace-Text "המאמר מציג מסגרת תיאורטית-אינפורמציונית חדשה לכימות מורכבות פיצ'רים במודלי דיפו ומציע גישה מתמטית להבנה פיצ'רים, מתי והיכן פיצ'רים מופיעים במהלך האימון. בניגוד לשיטות מסורתיות שמתמקדות בסליינסי (saliency) ושיוך פיצ'רים (attribution), המחקר מציע את מידת המורכבות שקיבלה שם v-information כמדד למורכבות חישובית, אשר מבטא את המאמץ הנדרש כדי לחלץ פיצ'רים במקום לשערך רק את התלות הסטטיסטית הישירה שלה בקלט." `
             ("המוטיבציה: " + $br + "דרישות ״החשמל״ של מערכות מבוססת AI, במיוחד מודלים גדולים, הופכות להיות יותר ויותר קשוחות. מכפלות נקודה צפה הן בין הפעולות החישוביות היקרות ביותר(מבחינת צריכת אנרגיה), והחלפתן באלטרנטיבות חסכוניות יותר יכולה להיות בעלת השלכות משמעותיות על תכנון חומרה למגוון רחב של יישומי AI. המחברים מדגישים כיצד צריכת האנרגיה ברשתות נוירונים עולה עם מספר פעולות הנקודה הצפה, ומכמתים את הפחתות האנרגיה האפשריות על ידי החלפת מכפלות בחיבורים.")

# --- Paragraph 4: technical basis of L-Mul ----------------------------------
Replace-Text "המחקר בוחן באופן שיטתי את התפתחותן בזמן אימון, התפלגותן המרחבית ותפקידן של פיצ'רים במודלים ויז'ן. הממצאים מצביעים על כך שמודלי דיפ מציגים תהליך למידה היררכי, שבו פיצ'רים פשוטים ודלות-מורכבות מופיעות מוקדם באימון ומתקדמות בקלות דרך חיבורים residual, בעוד פיצ'רים מורכבים יותר דורשות עיבוד עמוק יותר וזמן אימון ארוך יותר אך תורמות פחות משמעותית להחלטות הסופיות ממה שהיה מקובל להניח." `
             ("הבסיס הטכני של L-Mul:" + $br + " כפל נקודה צפה מסורתי כרוך בפעולות יקרות של מעריכי ומנטיסות. L-Mul עוקף זאת על ידי ארגון מחדש של החישוב, תוך שימוש בחיבור של מספרים שלמים במקום כפל של מנטיסות. המחברים תומכים בכך עם הערכת שגיאה תיאורטית, המראה ש-L-Mul עם מנטיסה של 3 ביטים מתעלה על מכפלת float8 e5m2, בעוד שעם מנטיסה של 4 ביטים הוא משתווה ואף מתעלה על float8 e4m3. דיוק מתמטי זה מספק אמינות חזקה לטענותיהם.")

# --- Paragraph 5: experiments -------------------------------------------------
Replace-Text "גישה מבוססת למורכבות בלמידת פיצ'רים" `
             ("ניסויים:" + $br + "המחברים משלבים את L-Mul בתוך מודלים מבוססי טרנספורמר ומעריכים את יעילותו במגוון משימות, כולל הבנת שפה טבעית, משימות הנמקה כלליות, ופתרון בעיות מתמטיות ועוד. יישום L-Mul למנגנון ה-attention מביא לאובדן דיוק זניח, ובמקרים מסוימים אף לשיפורים קלים בביצועים. המחברים אף מראים שהחלפת כל המכפלות בנקודה צפה במנטיסה של 3 ביטים בטרנספורמר מביאה לתוצאות דומות ל-float8 e4m3 הן בכיול (fine-tuning) והן בזמן הסקה.")

# --- Paragraph 6: pros and cons ------------------------------------------------
Replace-Text "ניתוח פיצ'רים בלמידה עמוקה התמקד עד כה בעיקר בחישוב החשיבות והשימושיות שלהן למשימה כזו או אחרת, אך כמעט ולא בוצע ניסיון לכמת כמה מורכב לחלץ פיצ'ר מתוך דאטה. מחקר זה משנה את נקודת המבט המסורתית בכך שהוא מציע מדד למאמץ חישובי הדרוש ללמידת פיצ'ר." `
             ("יתרונות וחסרונות:" + $br + "אחד ההיבטים המשכנעים ביותר במאמר הוא ההתמקדות ביעילות אנרגטית. על ידי שימוש בנתונים ממחקרים קודמים על צריכת אנרגיה בחומרה, המחברים מעריכים כי L-Mul יכול להפחית את עלות האנרגיה של מכפלות רכיביות ב-95% ואת עלות האנרגיה של פעולות מכפלה פנימית (dot product) ב-80%. זו טענה מרחיקת לכת, המציעה כי ל-L-Mul עשויות להיות השפעות מיידיות ומוחשיות על datacenters ויישומי AI בהיקפים גדולים.")

# --- Paragraph 7: open questions (no new w:br here, single run text) --------
Replace-Text "הגדרה מחדש של מורכבות פיצ'רים" `
             "המאמר מותיר כמה שאלות מעשיות ללא מענה. המחברים מכירים בכך של-GPUs קיימים אין תמיכה native ב-L-Mul, מה שמקשה על יישומו היעיל במערכות AI מודרניות. למרות שהמחברים רומזים כי חומרה ייעודית יכולה לאפשר אופטימיזציה של חישובי L-Mul, הם אינם מספקים תוכניות קונקרטיות פיתוחה."

# --- Paragraph 8: summary -----------------------------------------------------
Replace-Text "שיטות מסורתיות לשערוך פיצ'רים מסתמכות על שערוך מידע הדדי (mutual information) בין פיצ'ר לבין הדאטה. עם זאת, גישה זו אינה מביאה בחשבון את הקושי החישובי הכרוך בחילוץ הפיצ'ר." `
             ("סיכום:" + $br + "המאמר מציג גישה חדשנית להפחתת העלות החישובית והאנרגטית של LLMs ורשתות נוירונים אחרות. הביסוס התיאורטי חזק, תוצאות הניסוי משכנעות, וההשפעה הפוטנציאלית משמעותית. בעוד שנותרים אתגרים מעשיים—במיוחד באימוץ חומרה—עבודה זו פותחת דלתות חדשות לחישובי AI חסכוניים באנרגיה. אם תשופר ותאומץ, L-Mul עשוי למלא תפקיד מרכזי בהפיכת AI לבר-קיימא מבלי לפגוע בביצועים.")

# --- Remove the old paragraphs 9..33 (the rest of the original review body,
#     everything between the "summary" paragraph just rewritten above and the
#     final URL paragraph) -----------------------------------------------------
$firstToRemove = $d.Paragraphs.Item(9)
$lastToRemove  = $d.Paragraphs.Item(33)
$killRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
$killRange.Delete()

# --- Final paragraph: update the arxiv link ---------------------------------
Replace-Text "https://arxiv.org/abs/2407.06076" "https://arxiv.org/abs/2410.00907"

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
